$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 54742
$ws.Cells.Item(2, 5).Value = 708
$ws.Cells.Item(2, 6).Value = 708
$ws.Cells.Item(2, 7).Value = 1994
$ws.Cells.Item(2, 8).Value = -803
$ws.Cells.Item(2, 9).Value = -838
$ws.Cells.Item(2, 10).Value = 35
$ws.Cells.Item(2, 11).Value = 159686
$ws.Cells.Item(2, 12).Value = 41416
$ws.Cells.Item(2, 13).Value = 118269
$ws.Cells.Item(2, 14).Value = 115864
$ws.Cells.Item(2, 15).Value = 2405
$ws.Cells.Item(2, 16).Value = 3567
$ws.Cells.Item(2, 17).Value = 3106
$ws.Cells.Item(2, 18).Value = -3276
$ws.Cells.Item(2, 19).Value = -838
$ws.Cells.Item(2, 20).Value = 4763
$ws.Cells.Item(2, 21).Value = -1657
$ws.Cells.Item(2, 22).Value = 17778
$ws.Cells.Item(2, 23).Value = 1.29
$ws.Cells.Item(2, 24).Value = -1.47
$ws.Cells.Item(2, 25).Value = -0.88
$ws.Cells.Item(2, 26).Value = -0.61
$ws.Cells.Item(2, 27).Value = 35.02
$ws.Cells.Item(2, 28).Value = 2773.92
$ws.Cells.Item(2, 29).Value = -1426
$ws.Cells.Item(2, 30).Value = -81.36
$ws.Cells.Item(2, 31).Value = 164959
$ws.Cells.Item(2, 32).Value = 0.7
$ws.Cells.Item(2, 33).Value = 1000
$ws.Cells.Item(2, 34).Value = 0.86
$ws.Cells.Item(2, 35).Value = -83.86
$ws.Cells.Item(2, 36).Value = 68764530

# Row 3
$ws.Cells.Item(3, 4).Value = 49549
$ws.Cells.Item(3, 5).Value = -2675
$ws.Cells.Item(3, 6).Value = -598
$ws.Cells.Item(3, 7).Value = -1701
$ws.Cells.Item(3, 8).Value = 257
$ws.Cells.Item(3, 9).Value = 538
$ws.Cells.Item(3, 10).Value = -282
$ws.Cells.Item(3, 11).Value = 162253
$ws.Cells.Item(3, 12).Value = 49721
$ws.Cells.Item(3, 13).Value = 112532
$ws.Cells.Item(3, 14).Value = 110120
$ws.Cells.Item(3, 15).Value = 2412
$ws.Cells.Item(3, 16).Value = 3567
$ws.Cells.Item(3, 17).Value = 8811
$ws.Cells.Item(3, 18).Value = 1153
$ws.Cells.Item(3, 19).Value = -3547
$ws.Cells.Item(3, 20).Value = 7259
$ws.Cells.Item(3, 21).Value = 1552
$ws.Cells.Item(3, 22).Value = 17496
$ws.Cells.Item(3, 23).Value = -5.4
$ws.Cells.Item(3, 24).Value = 0.52
$ws.Cells.Item(3, 25).Value = 0.48
$ws.Cells.Item(3, 26).Value = 0.16
$ws.Cells.Item(3, 27).Value = 44.18
$ws.Cells.Item(3, 28).Value = 2770.97
$ws.Cells.Item(3, 29).Value = 765
$ws.Cells.Item(3, 30).Value = 149.01
$ws.Cells.Item(3, 31).Value = 156781
$ws.Cells.Item(3, 32).Value = 0.73
$ws.Cells.Item(3, 33).Value = 1000
$ws.Cells.Item(3, 34).Value = 0.88
$ws.Cells.Item(3, 35).Value = 130.58
$ws.Cells.Item(3, 36).Value = 68764530

# Row 4
$ws.Cells.Item(4, 4).Value = 52008
$ws.Cells.Item(4, 5).Value = -9263
$ws.Cells.Item(4, 6).Value = -9263
$ws.Cells.Item(4, 7).Value = -8207
$ws.Cells.Item(4, 8).Value = 2111
$ws.Cells.Item(4, 9).Value = 2194
$ws.Cells.Item(4, 10).Value = -83
$ws.Cells.Item(4, 11).Value = 149003
$ws.Cells.Item(4, 12).Value = 39362
$ws.Cells.Item(4, 13).Value = 109641
$ws.Cells.Item(4, 14).Value = 107221
$ws.Cells.Item(4, 15).Value = 2420
$ws.Cells.Item(4, 16).Value = 3567
$ws.Cells.Item(4, 17).Value = -13095
$ws.Cells.Item(4, 18).Value = 18543
$ws.Cells.Item(4, 19).Value = -8187
$ws.Cells.Item(4, 20).Value = 8326
$ws.Cells.Item(4, 21).Value = -21421
$ws.Cells.Item(4, 22).Value = 9505
$ws.Cells.Item(4, 23).Value = -17.81
$ws.Cells.Item(4, 24).Value = 4.06
$ws.Cells.Item(4, 25).Value = 2.02
$ws.Cells.Item(4, 26).Value = 1.36
$ws.Cells.Item(4, 27).Value = 35.9
$ws.Cells.Item(4, 28).Value = 2810.66
$ws.Cells.Item(4, 29).Value = 3117
$ws.Cells.Item(4, 30).Value = 34.97
$ws.Cells.Item(4, 31).Value = 158605
$ws.Cells.Item(4, 32).Value = 0.6899999999999999
$ws.Cells.Item(4, 33).Value = 1000
$ws.Cells.Item(4, 34).Value = 0.92
$ws.Cells.Item(4, 35).Value = 30.9
$ws.Cells.Item(4, 36).Value = 68764530

# Row 5
$ws.Cells.Item(5, 4).Value = 63466
$ws.Cells.Item(5, 5).Value = 1169
$ws.Cells.Item(5, 6).Value = 1169
$ws.Cells.Item(5, 7).Value = 8241
$ws.Cells.Item(5, 8).Value = 6432
$ws.Cells.Item(5, 9).Value = 6572
$ws.Cells.Item(5, 10).Value = -140
$ws.Cells.Item(5, 11).Value = 157417
$ws.Cells.Item(5, 12).Value = 42897
$ws.Cells.Item(5, 13).Value = 114520
$ws.Cells.Item(5, 14).Value = 112573
$ws.Cells.Item(5, 15).Value = 1947
$ws.Cells.Item(5, 16).Value = 3567
$ws.Cells.Item(5, 17).Value = -2501
$ws.Cells.Item(5, 18).Value = 893
$ws.Cells.Item(5, 19).Value = 3534
$ws.Cells.Item(5, 20).Value = 9915
$ws.Cells.Item(5, 21).Value = -12416
$ws.Cells.Item(5, 22).Value = 14246
$ws.Cells.Item(5, 23).Value = 1.84
$ws.Cells.Item(5, 24).Value = 10.13
$ws.Cells.Item(5, 25).Value = 5.98
$ws.Cells.Item(5, 26).Value = 4.2
$ws.Cells.Item(5, 27).Value = 37.46
$ws.Cells.Item(5, 28).Value = 2983.72
$ws.Cells.Item(5, 29).Value = 9338
$ws.Cells.Item(5, 30).Value = 21.9
$ws.Cells.Item(5, 31).Value = 168339
$ws.Cells.Item(5, 32).Value = 1.21
$ws.Cells.Item(5, 33).Value = 1000
$ws.Cells.Item(5, 34).Value = 0.49
$ws.Cells.Item(5, 35).Value = 10.19
$ws.Cells.Item(5, 36).Value = 68764530

# Row 6
$ws.Cells.Item(6, 4).Value = 91583
$ws.Cells.Item(6, 5).Value = 7150
$ws.Cells.Item(6, 6).Value = 7150
$ws.Cells.Item(6, 7).Value = 10362
$ws.Cells.Item(6, 8).Value = 7450
$ws.Cells.Item(6, 9).Value = 7012
$ws.Cells.Item(6, 11).Value = 193497
$ws.Cells.Item(6, 12).Value = 71245
$ws.Cells.Item(6, 13).Value = 122252
$ws.Cells.Item(6, 14).Value = 119340
$ws.Cells.Item(6, 16).Value = 3567
$ws.Cells.Item(6, 17).Value = 2606
$ws.Cells.Item(6, 18).Value = -17047
$ws.Cells.Item(6, 19).Value = 17561
$ws.Cells.Item(6, 20).Value = 21461
$ws.Cells.Item(6, 21).Value = -18855
$ws.Cells.Item(6, 22).Value = 32537
$ws.Cells.Item(6, 23).Value = 7.81
$ws.Cells.Item(6, 24).Value = 8.130000000000001
$ws.Cells.Item(6, 25).Value = 6.05
$ws.Cells.Item(6, 26).Value = 4.25
$ws.Cells.Item(6, 27).Value = 58.28
$ws.Cells.Item(6, 28).Value = 3266.07
$ws.Cells.Item(6, 29).Value = 9962
$ws.Cells.Item(6, 30).Value = 21.98
$ws.Cells.Item(6, 31).Value = 178459
$ws.Cells.Item(6, 32).Value = 1.23
$ws.Cells.Item(6, 33).Value = 1000
$ws.Cells.Item(6, 34).Value = 0.46
$ws.Cells.Item(6, 35).Value = 9.550000000000001
$ws.Cells.Item(6, 36).Value = 68764530

# Row 7
$ws.Cells.Item(7, 4).Value = 101655
$ws.Cells.Item(7, 5).Value = 4669
$ws.Cells.Item(7, 7).Value = 7188
$ws.Cells.Item(7, 8).Value = 5329
$ws.Cells.Item(7, 9).Value = 5028
$ws.Cells.Item(7, 11).Value = 206239
$ws.Cells.Item(7, 12).Value = 78755
$ws.Cells.Item(7, 13).Value = 127484
$ws.Cells.Item(7, 14).Value = 124352
$ws.Cells.Item(7, 16).Value = 3569
$ws.Cells.Item(7, 17).Value = 10721
$ws.Cells.Item(7, 18).Value = -15969
$ws.Cells.Item(7, 19).Value = 3791
$ws.Cells.Item(7, 20).Value = 20933
$ws.Cells.Item(7, 21).Value = -5528
$ws.Cells.Item(7, 23).Value = 4.59
$ws.Cells.Item(7, 24).Value = 5.24
$ws.Cells.Item(7, 25).Value = 4.13
$ws.Cells.Item(7, 26).Value = 2.67
$ws.Cells.Item(7, 27).Value = 61.78
$ws.Cells.Item(7, 29).Value = 7143
$ws.Cells.Item(7, 30).Value = 39.69
$ws.Cells.Item(7, 31).Value = 185954
$ws.Cells.Item(7, 32).Value = 1.52
$ws.Cells.Item(7, 33).Value = 989
$ws.Cells.Item(7, 34).Value = 0.35
$ws.Cells.Item(7, 35).Value = 13.52

# Row 8
$ws.Cells.Item(8, 4).Value = 118132
$ws.Cells.Item(8, 5).Value = 8614
$ws.Cells.Item(8, 7).Value = 11804
$ws.Cells.Item(8, 8).Value = 8975
$ws.Cells.Item(8, 9).Value = 8579
$ws.Cells.Item(8, 11).Value = 222459
$ws.Cells.Item(8, 12).Value = 87750
$ws.Cells.Item(8, 13).Value = 134709
$ws.Cells.Item(8, 14).Value = 131344
$ws.Cells.Item(8, 16).Value = 3569
$ws.Cells.Item(8, 17).Value = 13494
$ws.Cells.Item(8, 18).Value = -14513
$ws.Cells.Item(8, 19).Value = 4795
$ws.Cells.Item(8, 20).Value = 20011
$ws.Cells.Item(8, 21).Value = 96
$ws.Cells.Item(8, 23).Value = 7.29
$ws.Cells.Item(8, 24).Value = 7.6
$ws.Cells.Item(8, 25).Value = 6.71
$ws.Cells.Item(8, 26).Value = 4.19
$ws.Cells.Item(8, 27).Value = 65.14
$ws.Cells.Item(8, 29).Value = 12190
$ws.Cells.Item(8, 30).Value = 22.64
$ws.Cells.Item(8, 31).Value = 196409
$ws.Cells.Item(8, 32).Value = 1.41
$ws.Cells.Item(8, 33).Value = 1005
$ws.Cells.Item(8, 34).Value = 0.36
$ws.Cells.Item(8, 35).Value = 8.06

# Row 9
$ws.Cells.Item(9, 4).Value = 139376
$ws.Cells.Item(9, 5).Value = 11373
$ws.Cells.Item(9, 7).Value = 15285
$ws.Cells.Item(9, 8).Value = 11584
$ws.Cells.Item(9, 9).Value = 11119
$ws.Cells.Item(9, 11).Value = 241488
$ws.Cells.Item(9, 12).Value = 96163
$ws.Cells.Item(9, 13).Value = 145325
$ws.Cells.Item(9, 14).Value = 141773
$ws.Cells.Item(9, 16).Value = 3569
$ws.Cells.Item(9, 17).Value = 17666
$ws.Cells.Item(9, 18).Value = -15697
$ws.Cells.Item(9, 19).Value = 3781
$ws.Cells.Item(9, 20).Value = 20694
$ws.Cells.Item(9, 21).Value = 3758
$ws.Cells.Item(9, 23).Value = 8.16
$ws.Cells.Item(9, 24).Value = 8.31
$ws.Cells.Item(9, 25).Value = 8.140000000000001
$ws.Cells.Item(9, 26).Value = 4.99
$ws.Cells.Item(9, 27).Value = 66.17
$ws.Cells.Item(9, 29).Value = 15798
$ws.Cells.Item(9, 30).Value = 17.47
$ws.Cells.Item(9, 31).Value = 212005
$ws.Cells.Item(9, 33).Value = 1054
$ws.Cells.Item(9, 34).Value = 0.38
$ws.Cells.Item(9, 35).Value = 6.52
